# Applies the BGR model update (2025-08-15 22:44):
# Swap a handful of "cost class rank" values (column P) on the
# "solar" and "wind" sheets.

$wb = $excel.ActiveWorkbook

$wsSolar = $wb.Worksheets.Item("solar")
$wsWind  = $wb.Worksheets.Item("wind")

# --- solar sheet: swap P5 / P6 -----------------------------------------
$wsSolar.Range("P5").Value = 2
$wsSolar.Range("P6").Value = 4

# --- wind sheet: swap several pairs of values in column P --------------
$wsWind.Range("P4").Value  = 2
$wsWind.Range("P5").Value  = 3

$wsWind.Range("P15").Value = 1
$wsWind.Range("P17").Value = 3

$wsWind.Range("P18").Value = 2
$wsWind.Range("P19").Value = 1
$wsWind.Range("P20").Value = 3

$wsWind.Range("P47").Value = 1
$wsWind.Range("P48").Value = 2
